$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140-158 down to 141-159
$ws.Rows(140).Insert()

# Populate the newly inserted row 140 with the new record
$ws.Range("A140").Value2 = 10
$ws.Range("B140").Value = "Vega Modelo de Temuco"
$ws.Range("C140").Value = "La Araucanía"
$ws.Range("D140").Value2 = 44522
$ws.Range("E140").Value2 = 9
$ws.Range("F140").Value2 = 100112005
$ws.Range("G140").Value = "Puerro"
$ws.Range("H140").Value = "Azul de Maquehue"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value2 = 50
$ws.Range("K140").Value2 = 8000
$ws.Range("L140").Value2 = 8000
$ws.Range("M140").Value2 = 8000
$ws.Range("N140").Value = "$/docena de paquetes"
$ws.Range("O140").Value = "Provincia de Cautín"
$ws.Range("P140").Value2 = 667
$ws.Range("Q140").Value2 = 12
$ws.Range("R140").Value = "Hortaliza"
